$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(161).Insert()

$ws.Range("A161").Value = 5
$ws.Range("B161").Value = "Macroferia Regional de Talca"
$ws.Range("C161").Value = "Maule"
$ws.Range("D161").Value = 45089
$ws.Range("E161").Value = 7
$ws.Range("F161").Value = 100112031
$ws.Range("G161").Value = "Poroto verde"
$ws.Range("H161").Value = "Sin especificar"
$ws.Range("I161").Value = "Primera"
$ws.Range("J161").Value = 150
$ws.Range("K161").Value = 22000
$ws.Range("L161").Value = 22000
$ws.Range("M161").Value = 22000
$ws.Range("N161").Value = "$/malla 25 kilos"
$ws.Range("O161").Value = "Perú"
$ws.Range("P161").Value = 880
$ws.Range("Q161").Value = 25
$ws.Range("R161").Value = "Hortaliza"
